# ------------------------------------------------------------------
# DDAf_2023_Tableau_annexe_Tab03.xlsx - "Add files via upload" re-save
#
# The underlying source data for sheet "Tab03" was refreshed upstream
# and the file re-uploaded. The only semantically meaningful change
# (i.e. not a cosmetic artifact of a different Excel build re-saving
# the package - fileVersion/rupBuild, xr:uid/revisionPtr GUIDs, window
# geometry, calcPr, default font metrics driving dyDescent/row-height
# jitter) is a batch of updated numeric cell values on sheet "Tab03":
#   - rows 67, 68, 69, 70, 72, 73: last-digit recalculation jitter for
#     a handful of cells in the COMESA/CEN-SAD/CAE/CEEAC/IGAD/CDAA
#     regional-aggregate rows;
#   - rows 97-98 ("Afrique, Etats fragiles" / "RDM, Etats fragiles"):
#     the whole C:AP data row was recomputed with updated source data.
#
# All target cells hold literal pasted values (no formulas anywhere in
# the sheet), so the fix is a straight set of explicit cell writes.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab03")

# Map of row -> @{ column-letter = new value } taken from the updated
# workbook. Column letters are converted to 1-based indices below so
# the writes can go through $ws.Cells.Item(row, col).
$updates = @{
    67 = @{
        "F" = 2.5921470390226
        "G" = 2.38313691600849
        "K" = 2.53638052490155
        "L" = 2.59032532798567
        "M" = 2.61029847748309
        "N" = 2.63532591641213
        "O" = 2.66161119187305
        "R" = 2.66855184353583
        "S" = 2.68597831773045
        "U" = 2.72820966728888
        "W" = 2.72353652492812
        "AA" = 2.72687272279364
        "AB" = 2.67835112731019
        "AC" = 2.68966643927735
        "AE" = 2.6405015333437
        "AH" = 2.48727375040425
        "AI" = 2.42115407360322
        "AJ" = 2.40815488166162
        "AL" = 2.3812237511343
        "AM" = 2.35140904811455
        "AN" = 2.32382821362576
        "AO" = 2.29478300350006
    }
    68 = @{
        "C" = 2.51821184484806
        "D" = 2.34945487310925
        "G" = 2.49525319396742
        "H" = 2.51877297243319
        "I" = 2.47523024736915
        "J" = 2.43011280677512
        "Q" = 2.56469622516036
        "R" = 2.54542560918265
        "U" = 2.55060253829902
        "V" = 2.54355909651374
        "W" = 2.54988109904499
        "X" = 0.5023829046555
        "Y" = 2.41682304157733
        "Z" = 2.53559644557235
        "AA" = 2.55631722822514
        "AC" = 2.51279112342937
        "AH" = 2.31946344939391
        "AI" = 2.27418388666394
        "AK" = 2.26651269180713
        "AL" = 2.24909698093725
        "AN" = 2.19577840803864
    }
    69 = @{
        "C" = 3.08637079074761
        "K" = 2.63970756388729
        "L" = 2.78473637453445
        "U" = 3.06570855339148
        "AJ" = 2.78259415151685
        "AK" = 2.77392370596536
        "AM" = 2.72445842585771
        "AN" = 2.7014521015041
        "AO" = 2.67678583292925
    }
    70 = @{
        "E" = 3.07321114707839
        "F" = 2.70912694909911
        "G" = 2.10475783094328
        "H" = 2.54009690631938
        "I" = 3.32745985621559
        "J" = 2.96090899691044
        "R" = 3.23991727566111
        "S" = 3.23499552251581
        "T" = 3.25558047829184
        "W" = 3.38145325452473
        "Y" = 3.27851061762836
        "Z" = 3.28294784222918
        "AB" = 3.18024751265633
        "AC" = 3.19130683920168
        "AD" = 3.21284771440551
        "AF" = 3.11966859741137
        "AG" = 3.07677090015357
        "AJ" = 3.00533968628369
        "AM" = 2.9377982728513
        "AN" = 2.90949946191645
    }
    72 = @{
        "D" = 2.86165762708495
        "E" = 2.68481525314275
        "G" = 3.19162748339563
        "H" = 3.06767986155847
        "I" = 2.96962406407872
        "J" = 2.92163198639388
        "V" = 2.9239749896
        "W" = 2.8741537286102
        "AB" = 2.63283686179554
        "AC" = 2.6681517907986
        "AD" = 2.62457757552543
        "AF" = 2.67320388330869
        "AG" = 2.67703322663957
        "AI" = 2.52486399865648
    }
    73 = @{
        "C" = 2.99292939126716
        "F" = 2.60032722550592
        "H" = 2.83520380410929
        "K" = 2.23284493920597
        "S" = 2.48453182326027
        "U" = 2.55017627605243
        "AB" = 2.96628001543398
        "AC" = 2.81736179617322
        "AE" = 2.74355560951609
        "AF" = 2.69169190775176
        "AK" = 2.58974461795842
        "AL" = 2.58032788675349
        "AM" = 2.54930899732064
        "AN" = 2.52200532917484
        "AP" = 2.54708711239264
    }
    97 = @{
        "C" = 2.8712840883373
        "D" = 2.92492763925438
        "E" = 2.92669754493402
        "F" = 2.86790836050821
        "G" = 2.9526136339493
        "H" = 2.86883035789605
        "I" = 2.71223963431944
        "J" = 2.73054027080883
        "K" = 2.76854122895696
        "L" = 2.76301303564244
        "M" = 2.77832430535474
        "N" = 2.79846724662352
        "O" = 2.81189532078243
        "P" = 2.80383554361809
        "Q" = 2.79720548702116
        "R" = 2.80745484372817
        "S" = 2.82397736456947
        "T" = 2.83821503463311
        "U" = 2.83614176608713
        "V" = 2.8198287773505
        "W" = 2.83343994192984
        "X" = 2.79207742704781
        "Y" = 2.77947707258539
        "Z" = 2.81953564306341
        "AA" = 2.81405275856477
        "AB" = 2.81115965356371
        "AC" = 2.82615672421647
        "AD" = 2.82057692956375
        "AE" = 2.75508390380708
        "AF" = 2.70864799733255
        "AG" = 2.692456347053
        "AH" = 2.64443466147515
        "AI" = 2.59102143648207
        "AJ" = 2.57385748649384
        "AK" = 2.56067043243988
        "AL" = 2.53350918494584
        "AM" = 2.50188431638321
        "AN" = 2.47146444953703
        "AO" = 2.44038948714815
        "AP" = 2.50157463512106
    }
    98 = @{
        "C" = 2.75807477897359
        "D" = 2.66977998602422
        "E" = 2.25921752148226
        "F" = 2.01390369758676
        "G" = 2.04829064655117
        "H" = 2.14113598906958
        "I" = 2.10955296093132
        "J" = 2.07856312472203
        "K" = 2.0831190212121
        "L" = 2.07182583637118
        "M" = 2.13999660428088
        "N" = 2.1508015956939
        "O" = 1.84894495545929
        "P" = 1.7060571061198
        "Q" = 1.7815451200724
        "R" = 1.70927412964519
        "S" = 1.59230617603761
        "T" = 1.54736693405895
        "U" = 1.48749241447583
        "V" = 1.49593650226638
        "W" = 1.59036899341716
        "X" = 1.55037376729301
        "Y" = 1.45425020245993
        "Z" = 1.39423303169786
        "AA" = 1.3972830940858
        "AB" = 1.38546097810495
        "AC" = 1.2867654663596
        "AD" = 1.21234870234335
        "AE" = 1.14869949335288
        "AF" = 1.11090440242374
        "AG" = 1.1834309887613
        "AH" = 1.22828707688809
        "AI" = 1.28069165989329
        "AJ" = 1.35287377537074
        "AK" = 1.34477818678824
        "AL" = 1.31286237762476
        "AM" = 1.27934112749542
        "AN" = 1.24709579796447
        "AO" = 1.21512379955495
        "AP" = 1.2798298233897
    }
}

$colIndex = @{
    "C" = 3
    "D" = 4
    "E" = 5
    "F" = 6
    "G" = 7
    "H" = 8
    "I" = 9
    "J" = 10
    "K" = 11
    "L" = 12
    "M" = 13
    "N" = 14
    "O" = 15
    "P" = 16
    "Q" = 17
    "R" = 18
    "S" = 19
    "T" = 20
    "U" = 21
    "V" = 22
    "W" = 23
    "X" = 24
    "Y" = 25
    "Z" = 26
    "AA" = 27
    "AB" = 28
    "AC" = 29
    "AD" = 30
    "AE" = 31
    "AF" = 32
    "AG" = 33
    "AH" = 34
    "AI" = 35
    "AJ" = 36
    "AK" = 37
    "AL" = 38
    "AM" = 39
    "AN" = 40
    "AO" = 41
    "AP" = 42
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $colNum = $colIndex[$col]
        $ws.Cells.Item([int]$row, $colNum).Value = $updates[$row][$col]
    }
}
